$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles from column J into the new column K for the rows that need formatting,
# then set the new values (copy also carries the value, which we override afterwards
# where a different value is needed).
$ws.Range("J3").Copy($ws.Range("K3"))

$ws.Range("J4").Copy($ws.Range("K4"))
$ws.Range("K4").Value = 2021

$ws.Range("J5").Copy($ws.Range("K5"))
$ws.Range("K5").Value = 375

$ws.Range("J6").Copy($ws.Range("K6"))
$ws.Range("K6").Value = "-"

$ws.Range("J7").Copy($ws.Range("K7"))
$ws.Range("K7").Value = 5

$ws.Range("J8").Copy($ws.Range("K8"))
$ws.Range("K8").Value = "-"

$ws.Range("J9").Copy($ws.Range("K9"))
$ws.Range("K9").Value = 18

$ws.Range("J10").Copy($ws.Range("K10"))
$ws.Range("K10").Value = 150

# Reflect the new active cell selection on the sheet view (K7), matching the source edit.
$ws.Range("K7").Select()
